# Adds two new columns (C: tek_period_start_year, D: tek_period_end_year)
# to the TEK_parameters sheet, filling in the start/end year for each TEK
# period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (C1/D1).
$ws.Range("C1").Value = "tek_period_start_year"
$ws.Range("D1").Value = "tek_period_end_year"
# Header row (row 1) uses the same bold styling as the existing B1 header.
$ws.Range("C1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true

# Period start/end year values per TEK_ID row (rows 2-12).
$periods = @(
    @{ Row = 2;  Start = 0;    End = 1955 },
    @{ Row = 3;  Start = 1956; End = 1970 },
    @{ Row = 4;  Start = 1971; End = 1980 },
    @{ Row = 5;  Start = 1981; End = 1990 },
    @{ Row = 6;  Start = 1971; End = 1990 },
    @{ Row = 7;  Start = 1991; End = 2000 },
    @{ Row = 8;  Start = 2001; End = 2010 },
    @{ Row = 9;  Start = 2011; End = 2013 },
    @{ Row = 10; Start = 2014; End = 2020 },
    @{ Row = 11; Start = 2021; End = 2024 },
    @{ Row = 12; Start = 2025; End = 2040 }
)

foreach ($p in $periods) {
    $ws.Cells.Item($p.Row, 3).Value = $p.Start
    $ws.Cells.Item($p.Row, 4).Value = $p.End
}

# Match the column widths used for columns C/D in the saved workbook.
$ws.Columns.Item(3).ColumnWidth = 20
$ws.Columns.Item(4).ColumnWidth = 20

# Update print/page setup (A4, portrait) as captured in the saved workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the final selection on D12, matching the last edited cell.
$ws.Range("D12").Select()

$wb.Save()
